# Edit DNB Mastercard Demo workbook: rename sheet, widen date numFmt to
# include a time component, and replace the transaction rows with a more
# diverse merchant list (20 rows instead of 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet/tab from "DNB Mastercard Demo" to "Sheet".
$ws.Name = "Sheet"

# 2) Widen the custom date number format used by column A so it also
#    carries a time component (yyyy-mm-dd -> yyyy-mm-dd h:mm:ss).
$ws.Range("A2:A19").NumberFormat = "yyyy-mm-dd h:mm:ss"

# 3) Replace the transaction data (rows 2-14 become rows 2-19) with the
#    new, more diverse merchant list.
$data = @(
    @(45716, "SPOTIFY FAMILY", $null, 169),
    @(45715, "KIWI MINIPRIS LØKKA", $null, 289.5),
    @(45714, "NETFLIX PREMIUM", $null, 219),
    @(45713, "REMA TORSHOV", $null, 345),
    @(45712, "JULA OSLO", $null, 567),
    @(45710, "VINMONOPOLET AKER BRYGGE", $null, 489),
    @(45708, "STARBUCKS AKER BRYGGE", $null, 95),
    @(45706, "BILTEMA ALF", $null, 789),
    @(45705, "GITHUB.COM", $null, 129),
    @(45703, "COOP OBS BYGG ALNA", $null, 1456),
    @(45702, "IKEA FURUSET", $null, 3456),
    @(45700, "MENY TORSHOV", $null, 489.3),
    @(45698, "POWER ALNA", $null, 2199),
    @(45696, "KID INTERIØR STORO", $null, 899),
    @(45694, "CUBUS OSLO CITY", $null, 649),
    @(45693, "Innbetaling", 15000, $null),
    @(45691, "XXL SPORT CC VEST", $null, 1299),
    @(45689, "DRESSMANN MAJORSTUEN", $null, 799)
)

$row = 2
foreach ($entry in $data) {
    $date = $entry[0]
    $merchant = $entry[1]
    $inn = $entry[2]
    $ut = $entry[3]

    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $merchant

    if ($null -ne $inn) {
        $ws.Cells.Item($row, 5).Value = $inn
    } else {
        $ws.Cells.Item($row, 5).Value = $null
    }

    if ($null -ne $ut) {
        $ws.Cells.Item($row, 6).Value = $ut
    } else {
        $ws.Cells.Item($row, 6).Value = $null
    }

    $row++
}
